$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 8000.08
$ws.Range("D3").Value = 5000.0594

$ws.Range("D4").Value = 69075.08957500001
$ws.Range("E4").Value = "2025-03-27 15:32:22"

$ws.Range("D5").Value = 7000.0735
$ws.Range("E5").Value = "2025-03-27 17:21:11"
